$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 1502.8125
$ws.Range("I12").Value = 154.1
$ws.Range("J12").Value = 3750.6667
$ws.Range("K12").Value = 154.1
$ws.Range("L12").Value = 3750.6667
$ws.Range("M12").Value = 15.90000000000001
$ws.Range("N12").Value = -4090.6667
$ws.Range("H33").Value = 482.0435
$ws.Range("I33").Value = 513.9
$ws.Range("K33").Value = 513.9
$ws.Range("M33").Value = -284.9
$ws.Range("H76").Value = 90915810
$ws.Range("I76").Value = 6450
$ws.Range("J76").Value = 111117890
$ws.Range("K76").Value = 6450
$ws.Range("L76").Value = 111117890
$ws.Range("M76").Value = -6135
$ws.Range("N76").Value = -111118520
$ws.Range("H79").Value = 90915810
$ws.Range("I79").Value = 6450
$ws.Range("J79").Value = 111117890
$ws.Range("K79").Value = 6450
$ws.Range("L79").Value = 111117890
$ws.Range("M79").Value = -5358
$ws.Range("N79").Value = -111120074
$ws.Range("H98").Value = 2144.5454
$ws.Range("I98").Value = 1119.25
$ws.Range("J98").Value = 3374.9
$ws.Range("K98").Value = 1119.25
$ws.Range("L98").Value = 3374.9
$ws.Range("M98").Value = 378.75
$ws.Range("N98").Value = -6370.9
$ws.Range("H99").Value = 782.93335
$ws.Range("I99").Value = 288
$ws.Range("K99").Value = 864
$ws.Range("M99").Value = 634
$ws.Range("H122").Value = 2144.5454
$ws.Range("I122").Value = 1119.25
$ws.Range("J122").Value = 3374.9
$ws.Range("K122").Value = 3357.75
$ws.Range("L122").Value = 10124.7
$ws.Range("M122").Value = -907.75
$ws.Range("N122").Value = -15024.7
$ws.Range("H127").Value = 9988.333000000001
$ws.Range("I127").Value = 12890.667
$ws.Range("K127").Value = 38672.001
$ws.Range("M127").Value = -33712.001
$ws.Range("H132").Value = 16121.75
$ws.Range("I132").Value = 22597
$ws.Range("K132").Value = 67791
$ws.Range("M132").Value = -65261
$ws.Range("H135").Value = 2586.6
$ws.Range("J135").Value = 2983.25
$ws.Range("L135").Value = 26849.25
$ws.Range("N135").Value = -31919.25

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3164.6
$ws.Range("I32").Value = 3270.8684
$ws.Range("J32").Value = 1145.5
$ws.Range("K32").Value = 3270.8684
$ws.Range("L32").Value = 1145.5
$ws.Range("M32").Value = -2983.8684
$ws.Range("N32").Value = -1719.5
$ws.Range("H61").Value = 5094.0347
$ws.Range("I61").Value = 3431.1875
$ws.Range("J61").Value = 7140.615
$ws.Range("K61").Value = 3431.1875
$ws.Range("L61").Value = 7140.615
$ws.Range("M61").Value = -3219.1875
$ws.Range("N61").Value = -7564.615
$ws.Range("H97").Value = 998.53845
$ws.Range("I97").Value = 1061.2222
$ws.Range("K97").Value = 1061.2222
$ws.Range("M97").Value = -565.2221999999999
$ws.Range("H110").Value = 252357.84
$ws.Range("J110").Value = 2566
$ws.Range("L110").Value = 2566
$ws.Range("N110").Value = -6656
$ws.Range("H122").Value = 3328.0925
$ws.Range("I122").Value = 3899.389
$ws.Range("K122").Value = 11698.167
$ws.Range("M122").Value = -9248.167000000001
$ws.Range("H132").Value = 7914.8887
$ws.Range("I132").Value = 4031.0908
$ws.Range("J132").Value = 14018
$ws.Range("K132").Value = 12093.2724
$ws.Range("L132").Value = 42054
$ws.Range("M132").Value = -9563.2724
$ws.Range("N132").Value = -47114
$ws.Range("H136").Value = 5094.0347
$ws.Range("I136").Value = 3431.1875
$ws.Range("J136").Value = 7140.615
$ws.Range("K136").Value = 10293.5625
$ws.Range("L136").Value = 21421.845
$ws.Range("M136").Value = -7743.5625
$ws.Range("N136").Value = -26521.845

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H40").Value = 42000
$ws.Range("J40").Value = 42000
$ws.Range("L40").Value = 42000
$ws.Range("N40").Value = -42530
$ws.Range("H86").Value = 86914.586
$ws.Range("I86").Value = 145832.72
$ws.Range("J86").Value = 4429.2
$ws.Range("K86").Value = 145832.72
$ws.Range("L86").Value = 4429.2
$ws.Range("M86").Value = -144709.72
$ws.Range("N86").Value = -6675.2
$ws.Range("H89").Value = 86914.586
$ws.Range("I89").Value = 145832.72
$ws.Range("J89").Value = 4429.2
$ws.Range("K89").Value = 729163.6
$ws.Range("L89").Value = 22146
$ws.Range("M89").Value = -723547.6
$ws.Range("N89").Value = -33378
$ws.Range("H94").Value = 190.8
$ws.Range("I94").Value = 190.8
$ws.Range("K94").Value = 190.8
$ws.Range("M94").Value = 260.2
$ws.Range("H96").Value = 33475.668
$ws.Range("I96").Value = 33475.668
$ws.Range("K96").Value = 33475.668
$ws.Range("M96").Value = -30729.668
$ws.Range("H105").Value = 2907
$ws.Range("J105").Value = 4248.75
$ws.Range("L105").Value = 4248.75
$ws.Range("N105").Value = -7742.75
$ws.Range("H107").Value = 3812.5293
$ws.Range("I107").Value = 3602
$ws.Range("K107").Value = 3602
$ws.Range("M107").Value = -1682
$ws.Range("H134").Value = 67768.375
$ws.Range("I134").Value = 4411.875
$ws.Range("K134").Value = 13235.625
$ws.Range("M134").Value = -10700.625

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5779.4185
$ws.Range("I31").Value = 3029.1538
$ws.Range("J31").Value = 6971.2
$ws.Range("K31").Value = 3029.1538
$ws.Range("L31").Value = 6971.2
$ws.Range("M31").Value = -2734.1538
$ws.Range("N31").Value = -7561.2
$ws.Range("H34").Value = 5779.4185
$ws.Range("I34").Value = 3029.1538
$ws.Range("J34").Value = 6971.2
$ws.Range("K34").Value = 3029.1538
$ws.Range("L34").Value = 6971.2
$ws.Range("M34").Value = -2827.1538
$ws.Range("N34").Value = -7375.2
$ws.Range("H58").Value = 480993.2
$ws.Range("I58").Value = 1252182.9
$ws.Range("K58").Value = 1252182.9
$ws.Range("M58").Value = -1251979.9
$ws.Range("H94").Value = 694.6875
$ws.Range("I94").Value = 494.5
$ws.Range("J94").Value = 814.8
$ws.Range("K94").Value = 494.5
$ws.Range("L94").Value = 814.8
$ws.Range("M94").Value = -43.5
$ws.Range("N94").Value = -1716.8
$ws.Range("H99").Value = 7223.6665
$ws.Range("I99").Value = 6253.5
$ws.Range("J99").Value = 7999.8
$ws.Range("K99").Value = 6253.5
$ws.Range("L99").Value = 7999.8
$ws.Range("M99").Value = -4755.5
$ws.Range("N99").Value = -10995.8
$ws.Range("H126").Value = 7223.6665
$ws.Range("I126").Value = 6253.5
$ws.Range("J126").Value = 7999.8
$ws.Range("K126").Value = 18760.5
$ws.Range("L126").Value = 23999.4
$ws.Range("M126").Value = -16290.5
$ws.Range("N126").Value = -28939.4
$ws.Range("H132").Value = 8623.333000000001
$ws.Range("I132").Value = 9936.625
$ws.Range("J132").Value = 5996.75
$ws.Range("K132").Value = 29809.875
$ws.Range("L132").Value = 17990.25
$ws.Range("M132").Value = -27279.875
$ws.Range("N132").Value = -23050.25
$ws.Range("H136").Value = 480993.2
$ws.Range("I136").Value = 1252182.9
$ws.Range("K136").Value = 3756548.7
$ws.Range("M136").Value = -3753998.7

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 2898.2
$ws.Range("I3").Value = 2898.2
$ws.Range("K3").Value = 8694.599999999999
$ws.Range("M3").Value = -8582.599999999999
$ws.Range("H32").Value = 4275955.5
$ws.Range("I32").Value = 10000900
$ws.Range("J32").Value = 2367640.8
$ws.Range("K32").Value = 30002700
$ws.Range("L32").Value = 7102922.399999999
$ws.Range("M32").Value = -30002417
$ws.Range("N32").Value = -7103488.399999999
$ws.Range("H82").Value = 6875
$ws.Range("J82").Value = 7350
$ws.Range("L82").Value = 22050
$ws.Range("N82").Value = -22862
$ws.Range("H85").Value = 6875
$ws.Range("J85").Value = 7350
$ws.Range("L85").Value = 22050
$ws.Range("N85").Value = -24858
$ws.Range("H102").Value = 13901.1
$ws.Range("J102").Value = 14998.333
$ws.Range("L102").Value = 44994.999
$ws.Range("N102").Value = -49862.999
$ws.Range("H122").Value = 45467.39
$ws.Range("J122").Value = 57883.777
$ws.Range("L122").Value = 520953.993
$ws.Range("N122").Value = -525853.993
$ws.Range("H132").Value = 2946.4
$ws.Range("I132").Value = 1449.5
$ws.Range("J132").Value = 3944.3333
$ws.Range("K132").Value = 13045.5
$ws.Range("L132").Value = 35498.9997
$ws.Range("M132").Value = -10515.5
$ws.Range("N132").Value = -40558.9997
$ws.Range("H134").Value = 3299.4285
$ws.Range("I134").Value = 3299.4285
$ws.Range("K134").Value = 9898.2855
$ws.Range("M134").Value = -4828.2855
$ws.Range("H137").Value = 2168.077
$ws.Range("I137").Value = 2168.077
$ws.Range("J137").Value = 0
$ws.Range("K137").Value = 6504.231000000001
$ws.Range("L137").Value = 0
$ws.Range("M137").Value = -1404.231000000001
$ws.Range("N137").ClearContents()
$ws.Range("H140").Value = 3274.5
$ws.Range("I140").Value = 3274.5
$ws.Range("K140").Value = 9823.5
$ws.Range("M140").Value = -4643.5
$ws.Range("H141").Value = 2334.8333
$ws.Range("I141").Value = 2334.8333
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 7004.499899999999
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = -1824.499899999999
$ws.Range("N141").ClearContents()

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 488922.28
$ws.Range("I2").Value = 537764.6
$ws.Range("K2").Value = 537764.6
$ws.Range("M2").Value = -537651.6
$ws.Range("H70").Value = 5248.0386
$ws.Range("I70").Value = 5212.143
$ws.Range("K70").Value = 5212.143
$ws.Range("M70").Value = -4942.143
$ws.Range("H73").Value = 5248.0386
$ws.Range("I73").Value = 5212.143
$ws.Range("K73").Value = 5212.143
$ws.Range("M73").Value = -4276.143
$ws.Range("H80").Value = 784791.9
$ws.Range("I80").Value = 546661.5
$ws.Range("J80").Value = 1431145.6
$ws.Range("K80").Value = 546661.5
$ws.Range("L80").Value = 1431145.6
$ws.Range("M80").Value = -545663.5
$ws.Range("N80").Value = -1433141.6
$ws.Range("H83").Value = 784791.9
$ws.Range("I83").Value = 546661.5
$ws.Range("J83").Value = 1431145.6
$ws.Range("K83").Value = 2733307.5
$ws.Range("L83").Value = 7155728
$ws.Range("M83").Value = -2728315.5
$ws.Range("N83").Value = -7165712
$ws.Range("H102").Value = 4792.9
$ws.Range("I102").Value = 4833.6665
$ws.Range("K102").Value = 4833.6665
$ws.Range("M102").Value = -3211.6665
$ws.Range("H113").Value = 534803.4399999999
$ws.Range("I113").Value = 920424.9399999999
$ws.Range("K113").Value = 920424.9399999999
$ws.Range("M113").Value = -918254.9399999999
$ws.Range("H122").Value = 793484.1
$ws.Range("I122").Value = 923555.5600000001
$ws.Range("J122").Value = 13055.5
$ws.Range("K122").Value = 2770666.68
$ws.Range("L122").Value = 39166.5
$ws.Range("M122").Value = -2768216.68
$ws.Range("N122").Value = -44066.5
$ws.Range("H123").Value = 43999
$ws.Range("J123").Value = 43999
$ws.Range("L123").Value = 43999
$ws.Range("N123").Value = -48899
$ws.Range("H132").Value = 94711.27
$ws.Range("I132").Value = 3374.5
$ws.Range("J132").Value = 146903.72
$ws.Range("K132").Value = 10123.5
$ws.Range("L132").Value = 440711.16
$ws.Range("M132").Value = -7593.5
$ws.Range("N132").Value = -445771.16

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1200.2
$ws.Range("I16").Value = 1250.25
$ws.Range("J16").Value = 1000
$ws.Range("K16").Value = 1250.25
$ws.Range("L16").Value = 1000
$ws.Range("M16").Value = -1080.25
$ws.Range("N16").Value = -1340
$ws.Range("H40").Value = 7505126
$ws.Range("I40").Value = 12002802
$ws.Range("K40").Value = 12002802
$ws.Range("M40").Value = -12002666
$ws.Range("H46").Value = 3261.6897
$ws.Range("I46").Value = 3088.6843
$ws.Range("J46").Value = 3590.4
$ws.Range("K46").Value = 3088.6843
$ws.Range("L46").Value = 3590.4
$ws.Range("M46").Value = -2900.6843
$ws.Range("N46").Value = -3966.4
$ws.Range("H61").Value = 5100.88
$ws.Range("I61").Value = 3863.1333
$ws.Range("K61").Value = 3863.1333
$ws.Range("M61").Value = -3661.1333
$ws.Range("H102").Value = 0
$ws.Range("J102").Value = 0
$ws.Range("L102").Value = 0
$ws.Range("N102").ClearContents()
$ws.Range("H113").Value = 5100.88
$ws.Range("I113").Value = 3863.1333
$ws.Range("K113").Value = 3863.1333
$ws.Range("M113").Value = -1693.1333
$ws.Range("H132").Value = 7396.357
$ws.Range("I132").Value = 5412.25
$ws.Range("K132").Value = 16236.75
$ws.Range("M132").Value = -13706.75
$ws.Range("H136").Value = 6362.5
$ws.Range("I136").Value = 5000
$ws.Range("K136").Value = 15000
$ws.Range("M136").Value = -12450
$ws.Range("H138").Value = 74995
$ws.Range("J138").Value = 74995
$ws.Range("L138").Value = 74995
$ws.Range("N138").Value = -85275

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H69").Value = 20271
$ws.Range("J69").Value = 20271
$ws.Range("L69").Value = 20271
$ws.Range("N69").Value = -21769
$ws.Range("H72").Value = 20271
$ws.Range("J72").Value = 20271
$ws.Range("L72").Value = 60813
$ws.Range("N72").Value = -68301
$ws.Range("H75").Value = 33200
$ws.Range("J75").Value = 33200
$ws.Range("L75").Value = 33200
$ws.Range("N75").Value = -35072
$ws.Range("H78").Value = 33200
$ws.Range("J78").Value = 33200
$ws.Range("L78").Value = 99600
$ws.Range("N78").Value = -108960
$ws.Range("H81").Value = 18426.428
$ws.Range("I81").Value = 19690
$ws.Range("J81").Value = 2000
$ws.Range("K81").Value = 39380
$ws.Range("L81").Value = 4000
$ws.Range("M81").Value = -38319
$ws.Range("N81").Value = -6122
$ws.Range("H84").Value = 18426.428
$ws.Range("I84").Value = 19690
$ws.Range("J84").Value = 2000
$ws.Range("K84").Value = 196900
$ws.Range("L84").Value = 20000
$ws.Range("M84").Value = -191596
$ws.Range("N84").Value = -30608
$ws.Range("H98").Value = 266061.6
$ws.Range("J98").Value = 266061.6
$ws.Range("L98").Value = 266061.6
$ws.Range("N98").Value = -272051.6
$ws.Range("H113").Value = 669.73334
$ws.Range("I113").Value = 583.35486
$ws.Range("K113").Value = 1750.06458
$ws.Range("M113").Value = 419.9354199999998
$ws.Range("H122").Value = 52636620
$ws.Range("I122").Value = 111114580
$ws.Range("K122").Value = 333343740
$ws.Range("M122").Value = -333341290
$ws.Range("H126").Value = 2670.5625
$ws.Range("I126").Value = 1184.2106
$ws.Range("K126").Value = 3552.6318
$ws.Range("M126").Value = -1082.6318
$ws.Range("H132").Value = 15949.776
$ws.Range("I132").Value = 2001.1666
$ws.Range("J132").Value = 68257.06
$ws.Range("K132").Value = 6003.4998
$ws.Range("L132").Value = 204771.18
$ws.Range("M132").Value = -3473.4998
$ws.Range("N132").Value = -209831.18
$ws.Range("H136").Value = 327275.8
$ws.Range("I136").Value = 386881.22
$ws.Range("J136").Value = 186390.27
$ws.Range("K136").Value = 1160643.66
$ws.Range("L136").Value = 559170.8099999999
$ws.Range("M136").Value = -1158093.66
$ws.Range("N136").Value = -564270.8099999999

